$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently ends at row 20 (last weekly reading for this
# market/variety). A new weekly reading needs to be recorded, so the
# existing row 20 is pushed down to a new row 21, and row 20 is
# overwritten with the new reading's values.

# 1) Duplicate row 20 into row 21 (values + number format of the date
#    column), preserving the previous reading as a separate record.
for ($col = 1; $col -le 18; $col++) {
    $ws.Cells.Item(21, $col).Value2 = $ws.Cells.Item(20, $col).Value2
}
$ws.Cells.Item(21, 4).NumberFormat = $ws.Cells.Item(20, 4).NumberFormat

# 2) Overwrite row 20 with the new weekly reading.
$ws.Cells.Item(20, 4).Value2 = 44610    # Fecha
$ws.Cells.Item(20, 10).Value2 = 60      # Volumen
$ws.Cells.Item(20, 11).Value2 = 11000   # Precio minimo
$ws.Cells.Item(20, 12).Value2 = 12000   # Precio maximo
$ws.Cells.Item(20, 13).Value2 = 11500   # Precio promedio ponderado
$ws.Cells.Item(20, 15).Value2 = "Región Metropolitana"  # Origen
$ws.Cells.Item(20, 16).Value2 = 192     # Precio $/Kg
